$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 with same style (bold, centered, bordered) as the other headers
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (F2:F67)
$ws.Range("F2").Value = "2021-10-05 10:51:22.869695"
$ws.Range("F3").Value = "2021-10-05 10:51:22.869706"
$ws.Range("F4").Value = "2021-10-05 10:51:22.869710"
$ws.Range("F5").Value = "2021-10-05 10:51:22.869712"
$ws.Range("F6").Value = "2021-10-05 10:51:22.869715"
$ws.Range("F7").Value = "2021-10-05 10:51:22.869718"
$ws.Range("F8").Value = "2021-10-05 10:51:22.869720"
$ws.Range("F9").Value = "2021-10-05 10:51:22.869723"
$ws.Range("F10").Value = "2021-10-05 10:51:22.869725"
$ws.Range("F11").Value = "2021-10-05 10:51:22.869728"
$ws.Range("F12").Value = "2021-10-05 10:51:22.869730"
$ws.Range("F13").Value = "2021-10-05 10:51:22.869733"
$ws.Range("F14").Value = "2021-10-05 10:51:22.869735"
$ws.Range("F15").Value = "2021-10-05 10:51:22.869738"
$ws.Range("F16").Value = "2021-10-05 10:51:22.869740"
$ws.Range("F17").Value = "2021-10-05 10:51:22.869743"
$ws.Range("F18").Value = "2021-10-05 10:51:22.869746"
$ws.Range("F19").Value = "2021-10-05 10:51:22.869748"
$ws.Range("F20").Value = "2021-10-05 10:51:22.869751"
$ws.Range("F21").Value = "2021-10-05 10:51:22.869753"
$ws.Range("F22").Value = "2021-10-05 10:51:22.869756"
$ws.Range("F23").Value = "2021-10-05 10:51:22.869758"
$ws.Range("F24").Value = "2021-10-05 10:51:22.869761"
$ws.Range("F25").Value = "2021-10-05 10:51:22.869763"
$ws.Range("F26").Value = "2021-10-05 10:51:22.869766"
$ws.Range("F27").Value = "2021-10-05 10:51:22.869768"
$ws.Range("F28").Value = "2021-10-05 10:51:22.869771"
$ws.Range("F29").Value = "2021-10-05 10:51:22.869773"
$ws.Range("F30").Value = "2021-10-05 10:51:22.869776"
$ws.Range("F31").Value = "2021-10-05 10:51:22.869778"
$ws.Range("F32").Value = "2021-10-05 10:51:22.869781"
$ws.Range("F33").Value = "2021-10-05 10:51:22.869783"
$ws.Range("F34").Value = "2021-10-05 10:51:22.869786"
$ws.Range("F35").Value = "2021-10-05 10:51:22.869788"
$ws.Range("F36").Value = "2021-10-05 10:51:22.869791"
$ws.Range("F37").Value = "2021-10-05 10:51:22.869793"
$ws.Range("F38").Value = "2021-10-05 10:51:22.869796"
$ws.Range("F39").Value = "2021-10-05 10:51:22.869798"
$ws.Range("F40").Value = "2021-10-05 10:51:22.869801"
$ws.Range("F41").Value = "2021-10-05 10:51:22.869803"
$ws.Range("F42").Value = "2021-10-05 10:51:22.869806"
$ws.Range("F43").Value = "2021-10-05 10:51:22.869809"
$ws.Range("F44").Value = "2021-10-05 10:51:22.869811"
$ws.Range("F45").Value = "2021-10-05 10:51:22.869814"
$ws.Range("F46").Value = "2021-10-05 10:51:22.869816"
$ws.Range("F47").Value = "2021-10-05 10:51:22.869819"
$ws.Range("F48").Value = "2021-10-05 10:51:22.869821"
$ws.Range("F49").Value = "2021-10-05 10:51:22.869824"
$ws.Range("F50").Value = "2021-10-05 10:51:22.869826"
$ws.Range("F51").Value = "2021-10-05 10:51:22.869829"
$ws.Range("F52").Value = "2021-10-05 10:51:22.869831"
$ws.Range("F53").Value = "2021-10-05 10:51:22.869834"
$ws.Range("F54").Value = "2021-10-05 10:51:22.869837"
$ws.Range("F55").Value = "2021-10-05 10:51:22.869839"
$ws.Range("F56").Value = "2021-10-05 10:51:22.869842"
$ws.Range("F57").Value = "2021-10-05 10:51:22.869844"
$ws.Range("F58").Value = "2021-10-05 10:51:22.869847"
$ws.Range("F59").Value = "2021-10-05 10:51:22.869849"
$ws.Range("F60").Value = "2021-10-05 10:51:22.869852"
$ws.Range("F61").Value = "2021-10-05 10:51:22.869854"
$ws.Range("F62").Value = "2021-10-05 10:51:22.869857"
$ws.Range("F63").Value = "2021-10-05 10:51:22.869859"
$ws.Range("F64").Value = "2021-10-05 10:51:22.869862"
$ws.Range("F65").Value = "2021-10-05 10:51:22.869864"
$ws.Range("F66").Value = "2021-10-05 10:51:22.869868"
$ws.Range("F67").Value = "2021-10-05 10:51:22.869870"

$excel.CutCopyMode = 0

